# Grade update: scientific essays
# Adds a new "Science Paper" grade column (H) to the gradebook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("H1").Value = "Science Paper"

# New grade formulas (mirrors the existing "score/possible" style used in F & G)
$ws.Range("H2").Formula = "=0/100"
$ws.Range("H3").Formula = "=100/100"
$ws.Range("H4").Formula = "=0.9"
$ws.Range("H5").Formula = "=0.9"
$ws.Range("H8").Formula = "=0.93"
$ws.Range("H9").Formula = "=1"
$ws.Range("H11").Formula = "=95/100"

# Row 2 grew slightly taller once the new column was filled in
$ws.Rows.Item(2).RowHeight = 15

# Column sizing so the new header/content fits (closest attainable widths)
$ws.Columns.Item(1).ColumnWidth = 11.45
$ws.Columns.Item(2).ColumnWidth = 11.95
$ws.Columns.Item(8).ColumnWidth = 11.45

# Leave the selection where the last edit was made
$ws.Range("H10").Select() | Out-Null
